$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.432446002960205
$ws.Range("B1").Value = 3.867506265640259
$ws.Range("C1").Value = 2.921120882034302
$ws.Range("D1").Value = 2.777636289596558
$ws.Range("E1").Value = 2.496066093444824
